$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.577.36"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "1.664.96"
$ws.Range("E3").Value = "  -3.62%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'215.22"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").Value = "'0.513"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'23.51"
$ws.Range("E8").Value = "  -3.08%  "
$ws.Range("D9").Value = "'0.262"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "1.900.10"
$ws.Range("E12").Value = "  -3.62%  "
$ws.Range("D13").Value = "1.692.33"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("E14").Value = "  -2.77%  "
$ws.Range("D15").Value = "'0.556"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "'66.15"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("D17").Value = "'246.99"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").Value = "27.593.44"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "0.0₃0730"
$ws.Range("E19").Value = "  -3.73%  "
$ws.Range("D20").Value = "'7.53"
$ws.Range("E20").Value = "  -4.94%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  -3.78%  "
$ws.Range("E23").Value = "  -5.13%  "
$ws.Range("E24").Value = "  -5.19%  "
$ws.Range("D25").Value = "'145.97"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("E26").Value = "  -5.13%  "
$ws.Range("D27").Value = "'16.39"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "'0.111"
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("E30").Value = "  +3.62%  "
$ws.Range("D31").Value = "'0.0507"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("D33").Value = "1.473.76"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("D34").Value = "'3.10"
$ws.Range("E34").Value = "  -5.53%  "
$ws.Range("E35").Value = "  -5.88%  "
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0172"
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.573"
$ws.Range("E39").Value = "  -6.26%  "
$ws.Range("D40").Value = "'69.55"
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("E41").Value = "  -5.66%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'5.40"
$ws.Range("E43").Value = "  -7.61%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.21"
$ws.Range("E44").Value = "  -3.71%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.807.88"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("E47").Value = "  -4.26%  "
$ws.Range("D48").Value = "'89.14"
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("E49").Value = "  -2.35%  "
$ws.Range("E50").Value = "  -3.43%  "
$ws.Range("D51").Value = "'7.91"
$ws.Range("E51").Value = "  -4.01%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
